$d = $word.ActiveDocument

# --- Step 1: merge the "sub-groups" run and the bookmark out of the last
# paragraph, leaving one clean run with the full sentence. -------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$lastPara = $d.Paragraphs.Last
$fullRange = $lastPara.Range
$textOnly = $d.Range($fullRange.Start, $fullRange.End - 1)
$textOnly.Delete()
$insertPoint = $d.Range($fullRange.Start, $fullRange.Start)
$insertPoint.InsertAfter("Add articles regarding how to pool effect sizes across (multiple) sub-groups")

# --- Step 2: add the new list paragraph after it. ----------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.InsertAfter("Lav artikel om, hvordan du henter GitHub project ned")

# --- Step 3: re-create the _GoBack bookmark, collapsed, at the very end
# of the new paragraph's text (immediately before its paragraph mark).
# A collapsed Range built directly at that boundary confuses
# Bookmarks.Add, so instead: bookmark the final character (a valid,
# non-collapsed range), delete that character (which correctly collapses
# the bookmark to the boundary), then re-insert the character after the
# now-collapsed bookmark using a fresh Range (so the bookmark is left
# sitting right after the re-inserted text instead of wrapping it again).
$newPara = $d.Paragraphs.Last
$lastCharRange = $d.Range($newPara.Range.End - 2, $newPara.Range.End - 1)
$lastChar = $lastCharRange.Text
$d.Bookmarks.Add("_GoBack", $lastCharRange)
$bm = $d.Bookmarks("_GoBack")
$bm.Range.Delete()

$bm = $d.Bookmarks("_GoBack")
$reinsertPoint = $d.Range($bm.Start, $bm.Start)
$reinsertPoint.InsertAfter($lastChar)
